$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 955
$ws.Range("I28").Value = 573.5714
$ws.Range("J28").Value = 2290
$ws.Range("K28").Value = 573.5714
$ws.Range("L28").Value = 2290
$ws.Range("M28").Value = -88.57140000000004
$ws.Range("N28").Value = -3260
$ws.Range("H40").Value = 5004905
$ws.Range("I40").Value = 16671061
$ws.Range("J40").Value = 5124
$ws.Range("K40").Value = 16671061
$ws.Range("L40").Value = 5124
$ws.Range("M40").Value = -16670886
$ws.Range("N40").Value = -5474
$ws.Range("H99").Value = 854.7143
$ws.Range("I99").Value = 506
$ws.Range("J99").Value = 1319.6666
$ws.Range("K99").Value = 1518
$ws.Range("L99").Value = 3958.9998
$ws.Range("M99").Value = -20
$ws.Range("N99").Value = -6954.9998
$ws.Range("H101").Value = 930.1667
$ws.Range("I101").Value = 1073.6666
$ws.Range("J101").Value = 786.6667
$ws.Range("K101").Value = 3220.9998
$ws.Range("L101").Value = 2360.0001
$ws.Range("M101").Value = -1598.9998
$ws.Range("N101").Value = -5604.0001
$ws.Range("H107").Value = 686.4186
$ws.Range("I107").Value = 363.9
$ws.Range("J107").Value = 1430.6923
$ws.Range("K107").Value = 363.9
$ws.Range("L107").Value = 1430.6923
$ws.Range("M107").Value = 1556.1
$ws.Range("N107").Value = -5270.6923
$ws.Range("H110").Value = 43199.25
$ws.Range("J110").Value = 43199.25
$ws.Range("L110").Value = 43199.25
$ws.Range("N110").Value = -51379.25
$ws.Range("H132").Value = 3723.5952
$ws.Range("I132").Value = 3642.8462
$ws.Range("K132").Value = 10928.5386
$ws.Range("M132").Value = -8398.5386
$ws.Range("H138").Value = 317056.66
$ws.Range("J138").Value = 421997.28
$ws.Range("L138").Value = 1265991.84
$ws.Range("N138").Value = -1276271.84

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6997.8184
$ws.Range("J61").Value = 4990.1
$ws.Range("L61").Value = 4990.1
$ws.Range("N61").Value = -5414.1
$ws.Range("H74").Value = 4068.1428
$ws.Range("J74").Value = 1499
$ws.Range("L74").Value = 1499
$ws.Range("N74").Value = -3247
$ws.Range("H77").Value = 4068.1428
$ws.Range("J77").Value = 1499
$ws.Range("L77").Value = 7495
$ws.Range("N77").Value = -16231
$ws.Range("H113").Value = 49999
$ws.Range("J113").Value = 49999
$ws.Range("L113").Value = 49999
$ws.Range("N113").Value = -58677
$ws.Range("H122").Value = 2766.3914
$ws.Range("I122").Value = 2079.2144
$ws.Range("K122").Value = 6237.6432
$ws.Range("M122").Value = -3787.6432
$ws.Range("H136").Value = 6997.8184
$ws.Range("J136").Value = 4990.1
$ws.Range("L136").Value = 14970.3
$ws.Range("N136").Value = -20070.3

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 937.9259
$ws.Range("I20").Value = 749.8
$ws.Range("K20").Value = 749.8
$ws.Range("M20").Value = -502.8
$ws.Range("H94").Value = 3502.25
$ws.Range("I94").Value = 2779.6667
$ws.Range("K94").Value = 2779.6667
$ws.Range("M94").Value = -2328.6667
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 23093
$ws.Range("J28").Value = 23093
$ws.Range("L28").Value = 23093
$ws.Range("N28").Value = -23583
$ws.Range("H31").Value = 4798.077
$ws.Range("I31").Value = 6001.1113
$ws.Range("J31").Value = 4161.1763
$ws.Range("K31").Value = 6001.1113
$ws.Range("L31").Value = 4161.1763
$ws.Range("M31").Value = -5706.1113
$ws.Range("N31").Value = -4751.1763
$ws.Range("H34").Value = 4798.077
$ws.Range("I34").Value = 6001.1113
$ws.Range("J34").Value = 4161.1763
$ws.Range("K34").Value = 6001.1113
$ws.Range("L34").Value = 4161.1763
$ws.Range("M34").Value = -5799.1113
$ws.Range("N34").Value = -4565.1763
$ws.Range("H42").Value = 4999
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 4999
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 4999
$ws.Range("M42").ClearContents()
$ws.Range("N42").Value = -6185
$ws.Range("H58").Value = 5902.3687
$ws.Range("I58").Value = 4828.6875
$ws.Range("K58").Value = 4828.6875
$ws.Range("M58").Value = -4625.6875
$ws.Range("H134").Value = 4215.5
$ws.Range("I134").Value = 4231
$ws.Range("K134").Value = 12693
$ws.Range("M134").Value = -10158
$ws.Range("H136").Value = 5902.3687
$ws.Range("I136").Value = 4828.6875
$ws.Range("K136").Value = 14486.0625
$ws.Range("M136").Value = -11936.0625
$ws.Range("H141").Value = 123017.164
$ws.Range("J141").Value = 119939.4
$ws.Range("L141").Value = 119939.4
$ws.Range("N141").Value = -130299.4

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 2133.7222
$ws.Range("J107").Value = 2133.7222
$ws.Range("L107").Value = 6401.1666
$ws.Range("N107").Value = -10241.1666
$ws.Range("H113").Value = 522.7
$ws.Range("J113").Value = 623.75
$ws.Range("L113").Value = 1871.25
$ws.Range("N113").Value = -6211.25
$ws.Range("H132").Value = 1665.4615
$ws.Range("I132").Value = 1331
$ws.Range("J132").Value = 1814.1111
$ws.Range("K132").Value = 11979
$ws.Range("L132").Value = 16326.9999
$ws.Range("M132").Value = -9449
$ws.Range("N132").Value = -21386.9999

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5849.647
$ws.Range("I126").Value = 3183.4614
$ws.Range("J126").Value = 14514.75
$ws.Range("K126").Value = 9550.3842
$ws.Range("L126").Value = 43544.25
$ws.Range("M126").Value = -7080.3842
$ws.Range("N126").Value = -48484.25
$ws.Range("H132").Value = 5812.174
$ws.Range("I132").Value = 7485.625
$ws.Range("K132").Value = 22456.875
$ws.Range("M132").Value = -19926.875
$ws.Range("H134").Value = 29172
$ws.Range("J134").Value = 29172
$ws.Range("L134").Value = 87516
$ws.Range("N134").Value = -92586

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 45006
$ws.Range("I43").Value = 30012
$ws.Range("J43").Value = 60000
$ws.Range("K43").Value = 30012
$ws.Range("L43").Value = 60000
$ws.Range("M43").Value = -29819
$ws.Range("N43").Value = -60386
$ws.Range("H103").Value = 46804.25
$ws.Range("J103").Value = 46804.25
$ws.Range("L103").Value = 46804.25
$ws.Range("N103").Value = -49148.25

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 195536.19
$ws.Range("I2").Value = 16862.25
$ws.Range("J2").Value = 672000
$ws.Range("K2").Value = 16862.25
$ws.Range("L2").Value = 672000
$ws.Range("M2").Value = -16750.25
$ws.Range("N2").Value = -672224
$ws.Range("H4").Value = 8828.134
$ws.Range("I4").Value = 9906.546
$ws.Range("J4").Value = 5862.5
$ws.Range("K4").Value = 9906.546
$ws.Range("L4").Value = 5862.5
$ws.Range("M4").Value = -9793.546
$ws.Range("N4").Value = -6088.5
$ws.Range("H32").Value = 14605.2
$ws.Range("I32").Value = 14605.2
$ws.Range("K32").Value = 14605.2
$ws.Range("M32").Value = -14288.2
$ws.Range("H34").Value = 30026
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H81").Value = 59251.42
$ws.Range("I81").Value = 103528.3
$ws.Range("J81").Value = 10054.889
$ws.Range("K81").Value = 207056.6
$ws.Range("L81").Value = 20109.778
$ws.Range("M81").Value = -205995.6
$ws.Range("N81").Value = -22231.778
$ws.Range("H84").Value = 59251.42
$ws.Range("I84").Value = 103528.3
$ws.Range("J84").Value = 10054.889
$ws.Range("K84").Value = 1035283
$ws.Range("L84").Value = 100548.89
$ws.Range("M84").Value = -1029979
$ws.Range("N84").Value = -111156.89
$ws.Range("H96").Value = 2194.9167
$ws.Range("I96").Value = 1792.625
$ws.Range("K96").Value = 1792.625
$ws.Range("M96").Value = -419.625
$ws.Range("H126").Value = 2234.9
$ws.Range("I126").Value = 2280.6333
$ws.Range("J126").Value = 2097.7
$ws.Range("K126").Value = 6841.8999
$ws.Range("L126").Value = 6293.099999999999
$ws.Range("M126").Value = -4371.8999
$ws.Range("N126").Value = -11233.1
$ws.Range("H132").Value = 3962
$ws.Range("I132").Value = 4114.6665
$ws.Range("K132").Value = 12343.9995
$ws.Range("M132").Value = -9813.999500000002
$ws.Range("H136").Value = 3283.8518
$ws.Range("I136").Value = 1660.7333
$ws.Range("J136").Value = 5312.75
$ws.Range("K136").Value = 4982.199900000001
$ws.Range("L136").Value = 15938.25
$ws.Range("M136").Value = -2432.199900000001
$ws.Range("N136").Value = -21038.25
$ws.Range("H140").Value = 86660
$ws.Range("J140").Value = 86660
$ws.Range("L140").Value = 86660
$ws.Range("N140").Value = -97020
